$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.930.80"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.903.45"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.80"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4827"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3793"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07369"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9320"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.76"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.938.99"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.483"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.633"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008867"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "27.987.74"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.146"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "2.165.70"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.03"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.915"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.130"
$ws.Range("E28").Value = "  +6.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.25"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.964"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08936"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.263"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.254"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7661"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.667"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.544"
$ws.Range("E37").Value = "  -6.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.103"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5481"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05276"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.999"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.947"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1525"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.473"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.94"
$ws.Range("E45").Value = "  +6.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.646"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.89"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06082"
$ws.Range("E51").Value = "  -0.31%  "
